$wb = $excel.ActiveWorkbook

# --- Update status text "Ready for handoff" -> "In Translation" ---
# Overview sheet: columns E (zh-cn) and F (de-de), row 2
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

# zh-cn sheet: column C (Status), row 2
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

# de-de sheet: column C (Status), row 2
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# --- Narrow the "Status" columns ---
# The stored widths shrink from 17.2159881591797 to 13.4101845877511 characters.
# This runtime's ColumnWidth setter snaps to the nearest 1/6th of a character, so
# we pick the ColumnWidth input whose rounded result lands on the closest
# achievable width (13.333333333333334) to the target.
$newColumnWidth = 12.416666666666666

$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth

$wsZhCn.Columns.Item(3).ColumnWidth = $newColumnWidth

$wsDeDe.Columns.Item(3).ColumnWidth = $newColumnWidth
